$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 50

$ws.Cells.Item($row, 1).Value = "2024-01-12"
$ws.Cells.Item($row, 2).Value = "10:27:12"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "01"
$ws.Cells.Item($row, 5).Value = 139553
$ws.Cells.Item($row, 6).Value = 142843
$ws.Cells.Item($row, 7).Value = 171605
$ws.Cells.Item($row, 8).Value = 148287
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119367
$ws.Cells.Item($row, 11).Value = 224876
$ws.Cells.Item($row, 12).Value = 252501
$ws.Cells.Item($row, 13).Value = 185118
$ws.Cells.Item($row, 14).Value = 110556
$ws.Cells.Item($row, 15).Value = 40829
$ws.Cells.Item($row, 16).Value = 30898
$ws.Cells.Item($row, 17).Value = 72970
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42115
$ws.Cells.Item($row, 20).Value = -1
